$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 with new TPM-derived values ---
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.3971766666666667
$ws.Cells.Item(2, 8).Value = 1.19153

$ws.Cells.Item(2, 13).Value = 0.8377936666666667
$ws.Cells.Item(2, 14).Value = 2.513381
$ws.Cells.Item(2, 15).Value = 0.7130909380817101
$ws.Cells.Item(2, 16).Value = 0.7130909380817101
$ws.Cells.Item(2, 17).Value = 0.3327520958811111
$ws.Cells.Item(2, 18).Value = 2.99476886293
$ws.Cells.Item(2, 19).Value = 0.7130909380817101
$ws.Cells.Item(2, 20).Value = 0.7130909380817101

# --- Add new row 3 for the Resolving-Mac target cluster ---
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Pomc"
$ws.Cells.Item(3, 3).Value = "Oprm1"
$ws.Cells.Item(3, 4).Value = "Resolving-Mac"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.3971766666666667
$ws.Cells.Item(3, 8).Value = 1.19153
$ws.Cells.Item(3, 9).Value = 1
$ws.Cells.Item(3, 10).Value = 1
$ws.Cells.Item(3, 11).Value = 2
$ws.Cells.Item(3, 12).Value = 0.6666666666666666
$ws.Cells.Item(3, 13).Value = 0.3370826666666667
$ws.Cells.Item(3, 14).Value = 1.011248
$ws.Cells.Item(3, 15).Value = 0.2869090619182899
$ws.Cells.Item(3, 16).Value = 0.2869090619182899
$ws.Cells.Item(3, 17).Value = 0.1338813699377778
$ws.Cells.Item(3, 18).Value = 1.20493232944
$ws.Cells.Item(3, 19).Value = 0.2869090619182899
$ws.Cells.Item(3, 20).Value = 0.2869090619182899
